$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (refresh from GitHub Actions run)

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay as plain text like the source data.
$textCells = @("D4", "D5", "D6", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D31", "D32", "D34", "D35", "D37", "D39", "D42", "D43", "D45")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "55.089.53"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "2.289.73"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "505.72"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").Value = "129.06"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("D9").Value = "2.309.09"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +2.81%  "

$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("E12").Value = "  +8.25%  "

$ws.Range("E14").Value = "  +4.27%  "

$ws.Range("D15").Value = "2.719.82"
$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").Value = "55.228.00"
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("D18").Value = "2.287.94"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "312.92"
$ws.Range("E21").Value = "  +3.08%  "

$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  +4.48%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "60.16"
$ws.Range("E24").Value = "  -1.74%  "

$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "0.155"
$ws.Range("E26").Value = "  +3.08%  "

$ws.Range("E27").Value = "  +2.68%  "

$ws.Range("D28").Value = "172.32"
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("D29").Value = "0.0₃0711"
$ws.Range("E29").Value = "  +3.85%  "

$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  +3.73%  "

$ws.Range("D32").Value = "1.64"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "17.99"
$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("E36").Value = "  +3.46%  "

$ws.Range("D37").Value = "0.914"
$ws.Range("E37").Value = "  -4.78%  "

$ws.Range("E38").Value = "  +4.88%  "

$ws.Range("D39").Value = "36.87"
$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("E40").Value = "  +3.80%  "

$ws.Range("E41").Value = "  +1.13%  "

$ws.Range("D42").Value = "136.18"
$ws.Range("E42").Value = "  +9.03%  "

$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  +6.08%  "

$ws.Range("E44").Value = "  +2.21%  "

$ws.Range("D45").Value = "260.71"
$ws.Range("E45").Value = "  +9.23%  "

$ws.Range("E46").Value = "  +3.16%  "

$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("E50").Value = "  +3.27%  "

$ws.Range("E51").Value = "  +2.11%  "
